$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '68.204.17'
$ws.Range("E2").Value = '  -2.01%  '
$ws.Range("D3").Value = '3.580.57'
$ws.Range("E3").Value = '  -2.77%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '623.67'
$ws.Range("E5").Value = '  -6.50%  '
$ws.Range("D6").Value = '155.72'
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("D7").Value = '3.572.28'
$ws.Range("E7").Value = '  -2.89%  '
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -2.48%  '
$ws.Range("D11").Value = '6.96'
$ws.Range("E11").Value = '  -2.21%  '
$ws.Range("D12").Value = '0.434'
$ws.Range("E12").Value = '  -1.56%  '
$ws.Range("E13").Value = '  -2.76%  '
$ws.Range("D14").Value = '4.198.67'
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '32.24'
$ws.Range("E15").Value = '  -2.07%  '
$ws.Range("D16").Value = '3.595.81'
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '68.292.15'
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("E20").Value = '  -2.67%  '
$ws.Range("D21").Value = '459.95'
$ws.Range("E21").Value = '  -2.25%  '
$ws.Range("D22").Value = '9.84'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").Value = '0.645'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '78.09'
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").Value = '3.736.82'
$ws.Range("E25").Value = '  -2.42%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("D27").Value = '10.72'
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("E28").Value = '  -7.66%  '
$ws.Range("D29").Value = '8.42'
$ws.Range("E29").Value = '  -6.85%  '
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").Value = '1.64'
$ws.Range("E31").Value = '  -3.21%  '
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '26.11'
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("E34").Value = '  -4.28%  '
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.592.73'
$ws.Range("E35").Value = '  -2.50%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '6.22'
$ws.Range("E36").Value = '  -3.83%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.158'
$ws.Range("E37").Value = '  -4.84%  '
$ws.Range("D38").Value = '8.13'
$ws.Range("E38").Value = '  -3.88%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = '177.98'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("D42").Value = '5.66'
$ws.Range("E42").Value = '  -6.88%  '
$ws.Range("D43").Value = '0.0888'
$ws.Range("E43").Value = '  -2.15%  '
$ws.Range("E44").Value = '  -4.86%  '
$ws.Range("D45").Value = '0.901'
$ws.Range("E45").Value = '  -3.56%  '
$ws.Range("D46").Value = '46.04'
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("D47").Value = '28.43'
$ws.Range("E47").Value = '  +3.66%  '
$ws.Range("E48").Value = '  -4.98%  '
$ws.Range("D49").Value = '7.74'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  -5.52%  '
$ws.Range("E51").Value = '  -5.17%  '

$rng.Style = "Normal"
